# Risk management sheet update:
#  - Re-sort / re-enter a few data rows (rows 3 & 4 swap content, with a tweaked
#    probability value for "underestimated workload"; "loss of data" and
#    "loss of team members" probabilities bumped too).
#  - Shade rows 3, 5 and 7 white and row 4 light gray (manual banding) while
#    rows 2 and 6 stay unshaded.
#  - Row heights shrink on a couple of rows (less wrapped text now).
#  - A few columns widen slightly (bestFit after the longer text moved in).
#  - Selection changed to A1:G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-enter row data for rows 3, 4, 5, 7 (row 2 and row 6 keep their values).
# ---------------------------------------------------------------------------

# Row 3: "underestimated workload" (previously row 4), probability raised 0.4 -> 0.5
$ws.Range("A3").Value2 = "underestimated workload"
$ws.Range("B3").Value2 = "to much workload of this course regarding the workload of other courses"
$ws.Range("C3").Value2 = 0.5
$ws.Range("D3").Value2 = 5
$ws.Range("E3").Value2 = "flexible workload management"
$ws.Range("F3").Value2 = "Alex"
$ws.Range("G3").Formula = "=PRODUCT(Tabelle4[[#This Row],[prob of occurence (in %)]:[damage/impact (1-10)]])"

# Row 4: "server shutdown" (previously row 3), values unchanged
$ws.Range("A4").Value2 = "server shutdown"
$ws.Range("B4").Value2 = "shutdown due to failure/updates"
$ws.Range("C4").Value2 = 0.3
$ws.Range("D4").Value2 = 7
$ws.Range("E4").Value2 = "update freeze before endpresentation; backups; docker"
$ws.Range("F4").Value2 = "Mirko"
$ws.Range("G4").Formula = "=PRODUCT(Tabelle4[[#This Row],[prob of occurence (in %)]:[damage/impact (1-10)]])"

# Row 5: "loss of data" stays, probability raised 0.1 -> 0.2
$ws.Range("A5").Value2 = "loss of data "
$ws.Range("B5").Value2 = "due to device failure; "
$ws.Range("C5").Value2 = 0.2
$ws.Range("D5").Value2 = 10
$ws.Range("E5").Value2 = "code management using GitHub; track code sharing"
$ws.Range("F5").Value2 = "Alex"
$ws.Range("G5").Formula = "=PRODUCT(Tabelle4[[#This Row],[prob of occurence (in %)]:[damage/impact (1-10)]])"

# Row 7: "loss of team members" stays, probability raised 0.02 -> 0.07
$ws.Range("A7").Value2 = "loss of team members"
$ws.Range("B7").Value2 = "loss of team members due to exmatriculation"
$ws.Range("C7").Value2 = 0.07
$ws.Range("D7").Value2 = 8
$ws.Range("E7").Value2 = "share knowledge via GitHub; document detailed"
$ws.Range("F7").Value2 = "Alex"
$ws.Range("G7").Formula = "=PRODUCT(Tabelle4[[#This Row],[prob of occurence (in %)]:[damage/impact (1-10)]])"

# ---------------------------------------------------------------------------
# 2) Row heights (rows 2 and 6 get shorter; 3/4/5/7 reflect the re-ordered text)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 57
$ws.Rows.Item(3).RowHeight = 42.75
$ws.Rows.Item(4).RowHeight = 42.75
$ws.Rows.Item(5).RowHeight = 28.5
$ws.Rows.Item(6).RowHeight = 57
$ws.Rows.Item(7).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 3) Column widths (bestFit-style widen on B, D, E, G)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 31.666666666666668
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws.Columns.Item(5).ColumnWidth = 22
$ws.Columns.Item(7).ColumnWidth = 13.333333333333334

# ---------------------------------------------------------------------------
# 4) Manual row banding: rows 3, 5, 7 -> white fill; row 4 -> light gray fill
# ---------------------------------------------------------------------------
$ws.Range("A3:G3").Interior.Color = 16777215
$ws.Range("A5:G5").Interior.Color = 16777215
$ws.Range("A7:G7").Interior.Color = 16777215
$ws.Range("A4:G4").Interior.Color = 14277081

# ---------------------------------------------------------------------------
# 5) Selection
# ---------------------------------------------------------------------------
$ws.Range("A1:G7").Select()
